$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1577.1154
$ws.Range("I19").Value = 936
$ws.Range("J19").Value = 2325.0833
$ws.Range("K19").Value = 936
$ws.Range("L19").Value = 2325.0833
$ws.Range("M19").Value = -761
$ws.Range("N19").Value = -2675.0833

$ws.Range("H86").Value = 7719.5
$ws.Range("I86").Value = 6565.1665
$ws.Range("J86").Value = 9451
$ws.Range("K86").Value = 6565.1665
$ws.Range("L86").Value = 9451
$ws.Range("M86").Value = -5442.1665
$ws.Range("N86").Value = -11697

$ws.Range("H89").Value = 7719.5
$ws.Range("I89").Value = 6565.1665
$ws.Range("J89").Value = 9451
$ws.Range("K89").Value = 32825.8325
$ws.Range("L89").Value = 47255
$ws.Range("M89").Value = -27209.8325
$ws.Range("N89").Value = -58487

$ws.Range("H100").Value = 2875
$ws.Range("I100").Value = 3000
$ws.Range("J100").Value = 2833.3333
$ws.Range("K100").Value = 3000
$ws.Range("L100").Value = 2833.3333
$ws.Range("M100").Value = -2459
$ws.Range("N100").Value = -3915.3333

$ws.Range("H106").Value = 5933.864
$ws.Range("I106").Value = 953.75
$ws.Range("J106").Value = 7040.5557
$ws.Range("K106").Value = 953.75
$ws.Range("L106").Value = 7040.5557
$ws.Range("M106").Value = -322.75
$ws.Range("N106").Value = -8302.555700000001

$ws.Range("H141").Value = 1054.6786
$ws.Range("I141").Value = 908.55554
$ws.Range("J141").Value = 5000
$ws.Range("K141").Value = 2725.66662
$ws.Range("L141").Value = 15000
$ws.Range("M141").Value = 2454.33338
$ws.Range("N141").Value = -25360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2590.52
$ws.Range("I32").Value = 2027.9121
$ws.Range("J32").Value = 8279.111000000001
$ws.Range("K32").Value = 2027.9121
$ws.Range("L32").Value = 8279.111000000001
$ws.Range("M32").Value = -1740.9121
$ws.Range("N32").Value = -8853.111000000001

$ws.Range("H61").Value = 1194.2667
$ws.Range("I61").Value = 1035.193
$ws.Range("J61").Value = 1698
$ws.Range("K61").Value = 1035.193
$ws.Range("L61").Value = 1698
$ws.Range("M61").Value = -823.193
$ws.Range("N61").Value = -2122

$ws.Range("H88").Value = 1610.4615
$ws.Range("I88").Value = 904.5
$ws.Range("J88").Value = 2740
$ws.Range("K88").Value = 904.5
$ws.Range("L88").Value = 2740
$ws.Range("M88").Value = -498.5
$ws.Range("N88").Value = -3552

$ws.Range("H91").Value = 1610.4615
$ws.Range("I91").Value = 904.5
$ws.Range("J91").Value = 2740
$ws.Range("K91").Value = 904.5
$ws.Range("L91").Value = 2740
$ws.Range("M91").Value = 499.5
$ws.Range("N91").Value = -5548

$ws.Range("H136").Value = 1194.2667
$ws.Range("I136").Value = 1035.193
$ws.Range("J136").Value = 1698
$ws.Range("K136").Value = 3105.579
$ws.Range("L136").Value = 5094
$ws.Range("M136").Value = -555.5789999999997
$ws.Range("N136").Value = -10194

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H23").Value = 33200
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 33200
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 33200
$ws.Range("N23").Value = -33766

$ws.Range("H134").Value = 536213.25
$ws.Range("I134").Value = 891718.75
$ws.Range("J134").Value = 2955
$ws.Range("K134").Value = 2675156.25
$ws.Range("L134").Value = 8865
$ws.Range("M134").Value = -2672621.25
$ws.Range("N134").Value = -13935

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9435689
$ws.Range("I31").Value = 1147.5883
$ws.Range("J31").Value = 26318554
$ws.Range("K31").Value = 1147.5883
$ws.Range("L31").Value = 26318554
$ws.Range("M31").Value = -852.5882999999999
$ws.Range("N31").Value = -26319144

$ws.Range("H34").Value = 9435689
$ws.Range("I34").Value = 1147.5883
$ws.Range("J34").Value = 26318554
$ws.Range("K34").Value = 1147.5883
$ws.Range("L34").Value = 26318554
$ws.Range("M34").Value = -945.5882999999999
$ws.Range("N34").Value = -26318958

$ws.Range("H141").Value = 88500
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 88500
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 88500
$ws.Range("N141").Value = -98860

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 249
$ws.Range("I17").Value = 198.66667
$ws.Range("J17").Value = 400
$ws.Range("K17").Value = 596.00001
$ws.Range("L17").Value = 1200
$ws.Range("M17").Value = -427.00001
$ws.Range("N17").Value = -1538

$ws.Range("H46").Value = 751902.8
$ws.Range("I46").Value = 561.2
$ws.Range("J46").Value = 1002350
$ws.Range("K46").Value = 1683.6
$ws.Range("L46").Value = 3007050
$ws.Range("M46").Value = -1592.6
$ws.Range("N46").Value = -3007232

$ws.Range("H61").Value = 391.7
$ws.Range("I61").Value = 158.22223
$ws.Range("J61").Value = 582.7273
$ws.Range("K61").Value = 474.66669
$ws.Range("L61").Value = 1748.1819
$ws.Range("M61").Value = -259.66669
$ws.Range("N61").Value = -2178.1819

$ws.Range("H114").Value = 11231433
$ws.Range("I114").Value = 11111445
$ws.Range("J114").Value = 11366419
$ws.Range("K114").Value = 33334335
$ws.Range("L114").Value = 34099257
$ws.Range("M114").Value = -33331081
$ws.Range("N114").Value = -34105765

$ws.Range("H117").Value = 794.75
$ws.Range("I117").Value = 590
$ws.Range("J117").Value = 999.5
$ws.Range("K117").Value = 1770
$ws.Range("L117").Value = 2998.5
$ws.Range("M117").Value = 1672
$ws.Range("N117").Value = -9882.5

$ws.Range("H121").Value = 1691715.4
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 1691715.4
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 5075146.199999999
$ws.Range("M121").ClearContents()
$ws.Range("N121").Value = -5077766.199999999

$ws.Range("H129").Value = 3032
$ws.Range("I129").Value = 4685
$ws.Range("J129").Value = 2205.5
$ws.Range("K129").Value = 14055
$ws.Range("L129").Value = 6616.5
$ws.Range("M129").Value = -9055
$ws.Range("N129").Value = -16616.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 5320.8184
$ws.Range("I132").Value = 5789.5264
$ws.Range("J132").Value = 4684.7144
$ws.Range("K132").Value = 17368.5792
$ws.Range("L132").Value = 14054.1432
$ws.Range("M132").Value = -14838.5792
$ws.Range("N132").Value = -19114.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1944.3334
$ws.Range("I81").Value = 933.1667
$ws.Range("J81").Value = 3966.6667
$ws.Range("K81").Value = 1866.3334
$ws.Range("L81").Value = 7933.3334
$ws.Range("M81").Value = -805.3334
$ws.Range("N81").Value = -10055.3334

$ws.Range("H84").Value = 1944.3334
$ws.Range("I84").Value = 933.1667
$ws.Range("J84").Value = 3966.6667
$ws.Range("K84").Value = 9331.666999999999
$ws.Range("L84").Value = 39666.667
$ws.Range("M84").Value = -4027.666999999999
$ws.Range("N84").Value = -50274.667
